$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 24 for Cameroon, shifting the
#        existing rows 24-28 down to 25-29 ---
$ws.Rows.Item(24).Insert()

$newRowValues = @("'1","'1cLaLs","'1","'23","'2024","'39","Cameroon","CMR","CMR","'134","Mauritius","MTS","MUS","'0","'5","'0","'0","'0","'0","'0","-","'0")
for ($i = 0; $i -lt $newRowValues.Length; $i++) {
    $ws.Cells.Item(24, $i + 1).Value2 = $newRowValues[$i]
}

# --- 2. Fix up the item-sequence numbers (col D) for the rows that were
#        pushed down by the insert ---
$ws.Range("D25").Value2 = "'24"
$ws.Range("D26").Value2 = "'25"
$ws.Range("D27").Value2 = "'26"
$ws.Range("D28").Value2 = "'27"
$ws.Range("D29").Value2 = "'28"

# --- 3. Data updates ---
# Burundi row: asylum_seekers (col O) 5 -> 20
$ws.Range("O22").Value2 = "'20"
# Dem. Rep. of the Congo row (now row 25): refugees (col N) 0 -> 5,
# asylum_seekers (col O) 13 -> 27
$ws.Range("N25").Value2 = "'5"
$ws.Range("O25").Value2 = "'27"
# Nigeria row (now row 27): asylum_seekers (col O) 10 -> 5
$ws.Range("O27").Value2 = "'5"

# --- 4. Refresh the "short-url" value (col B) across every data row (the
#        shared string used to read "1cLaLs", now "prA9X8") ---
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cur = $ws.Cells.Item($r, 2).Value2
    if ($cur -eq "1cLaLs") {
        $ws.Cells.Item($r, 2).Value2 = "'prA9X8"
    }
}

# --- 5. Fix-up pass: re-apply the sheet's normal per-column number
#        format/style (using an untouched row as the formatting source)
#        to every row touched above, since forcing text entry above can
#        otherwise leave a stray "text" number format behind ---
$ws.Range("A2:V2").Copy()
$ws.Range("A22:V29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
